$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.093.89"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.90%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.583.60"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.90%  "

$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.37"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.06"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.56%  "

$ws.Range("E8").Value = "  -5.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.588.37"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.31%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.78"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +8.49%  "

$ws.Range("E11").Value = "  -1.45%  "

$ws.Range("E12").Value = "  +0.75%  "

$ws.Range("E13").Value = "  +0.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.034.28"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.109.36"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.96%  "

$ws.Range("E16").Value = "  -2.25%  "

$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.582.31"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.60%  "

$ws.Range("E19").Value = "  +1.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.64"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.31"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.06"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("E23").Value = "  -0.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.15"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.86%  "

$ws.Range("E25").Value = "  +0.24%  "

$ws.Range("E26").Value = "  -1.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.686.14"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.52%  "

$ws.Range("E29").Value = "  +2.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.45"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.56%  "

$ws.Range("E31").Value = "  +0.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.02"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.25"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.10%  "

$ws.Range("E34").Value = "  -1.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.73"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.18"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.863"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +20.20%  "

$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.842"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.89%  "

$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.77"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "300.05"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.03%  "

$ws.Range("E42").Value = "  +0.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "35.58"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.54%  "

$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.617"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.30%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0996"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.63%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0561"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.998"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.71"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.09%  "

$ws.Range("E49").Value = "  +1.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0233"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.012.41"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.73%  "
